# Generate Report for Handoff
# - Flip the localization status from "In Translation" to "Ready for handoff"
#   on every sheet that surfaces it (Overview status columns for zh-cn/de-de,
#   plus the per-language "Status" column on the zh-cn and de-de sheets).
# - Refresh the handoff timestamps that accompany that status change.
# - Widen the status columns so the new, longer "Ready for handoff" text
#   still fits (mirrors the column autosize that happens when the report is
#   regenerated).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---------------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status column
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status column
$wsZhCn.Range("C2").Value     = "Ready for handoff"   # Status column
$wsDeDe.Range("C2").Value     = "Ready for handoff"   # Status column

# --- Refresh the "Latest Handoff" timestamps ------------------------------
$wsOverview.Range("G2").Value = "2016-09-04 20:45:25"  # Latest HO Xliff Generate Date
$wsDeDe.Range("H2").Value     = "2016-09-04 20:45:25"  # Latest Handoff Datetime (de-de)
$wsZhCn.Range("H2").Value     = "2016-09-04 20:45:21"  # Latest Handoff Datetime (zh-cn)

# --- Widen the status columns to fit "Ready for handoff" ------------------
# Target display width ~17.22 characters. Excel's ColumnWidth setter snaps
# to whole-pixel increments, so we dial in the input that lands on the
# closest achievable grid value (~17.17) to the desired width.
$targetStatusColumnWidth = 16.33
$wsOverview.Columns.Item(5).ColumnWidth = $targetStatusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetStatusColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth     = $targetStatusColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth     = $targetStatusColumnWidth
